# Updated symbol list on Wed Feb  1 10:40:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.79"
$ws.Range("E2").Value = "'-1.21%"

$ws.Range("D3").Value = "'36.92"
$ws.Range("E3").Value = "'-0.64%"

$ws.Range("D4").Value = "'5.114"
$ws.Range("E4").Value = "'-0.22%"

$ws.Range("D5").Value = "'0.07799"
$ws.Range("E5").Value = "'0.45%"

$ws.Range("D6").Value = "'8.227"
$ws.Range("E6").Value = "'0.15%"

$ws.Range("D7").Value = "'1.877"
$ws.Range("E7").Value = "'0.21%"

$ws.Range("E8").Value = "'-2.34%"

$ws.Range("D9").Value = "'0.9231"
$ws.Range("E9").Value = "'0.29%"

$ws.Range("D10").Value = "'0.1109"
$ws.Range("E10").Value = "'-8.40%"

$ws.Range("E11").Value = "'-0.13%"

$ws.Range("D12").Value = "'0.08832"
$ws.Range("E12").Value = "'-4.22%"

$ws.Range("D13").Value = "'0.03356"
$ws.Range("E13").Value = "'-2.12%"

$ws.Range("E14").Value = "'-0.96%"

$ws.Range("D15").Value = "'0.001387"
$ws.Range("E15").Value = "'1.65%"

$ws.Range("D16").Value = "'0.005763"
$ws.Range("E16").Value = "'-3.02%"

$ws.Range("D17").Value = "'3.438"
$ws.Range("E17").Value = "'-3.42%"

$ws.Range("D18").Value = "'4.397"

$ws.Range("D19").Value = "'0.3422"

$ws.Range("D20").Value = "'6.243"
$ws.Range("E20").Value = "'18.70%"

$ws.Range("D21").Value = "'0.1281"
$ws.Range("E21").Value = "'0.97%"

$ws.Range("E22").Value = "'-6.03%"

$ws.Range("D23").Value = "'0.04323"
$ws.Range("E23").Value = "'-1.20%"

$ws.Range("E24").Value = "'-0.15%"

$ws.Range("D25").Value = "'0.004242"
$ws.Range("E25").Value = "'-0.24%"

$ws.Range("E26").Value = "'0.03%"

$ws.Range("D39").Value = "'0.02148"
$ws.Range("E39").Value = "'2.60%"

$ws.Range("D40").Value = "'0.05009"
$ws.Range("E40").Value = "'-0.35%"

$ws.Range("D41").Value = "'0.007551"
$ws.Range("E41").Value = "'-1.94%"

$ws.Range("D42").Value = "'0.1345"
$ws.Range("E42").Value = "'-0.19%"

$ws.Range("E43").Value = "'-11.46%"

$ws.Range("D44").Value = "'0.002053"
$ws.Range("E44").Value = "'-5.51%"

$ws.Range("D45").Value = "'0.008763"
$ws.Range("E45").Value = "'-8.45%"

$ws.Range("D46").Value = "'0.00006563"
$ws.Range("E46").Value = "'-2.24%"

$ws.Range("E47").Value = "'0.08%"

$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002828"
$ws.Range("E48").Value = "'-3.72%"

$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001002"
$ws.Range("E49").Value = "'-16.55%"

$ws.Range("E50").Value = "'0.08%"

$ws.Range("E51").Value = "'0.08%"
